$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New "Variant" header block: L3 (label) + M3:P3 (variant numbers 1..4)
# ---------------------------------------------------------------------------
$ws.Range("L3").Value = "Вариант"
$ws.Range("L3").Font.Bold = $true
$ws.Range("L3").HorizontalAlignment = -4108
$ws.Range("L3").VerticalAlignment = -4108
$ws.Range("L3").WrapText = $true

$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 2
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 4

# ---------------------------------------------------------------------------
# 2) Row 4 totals per variant
# ---------------------------------------------------------------------------
$ws.Range("M4").Formula = "=SUM(M5:M33)"
$ws.Range("N4").Formula = "=SUM(N5:N33)"
$ws.Range("O4").Formula = "=SUM(O5:O33)"
$ws.Range("P4").Formula = "=SUM(P5:P33)"

# ---------------------------------------------------------------------------
# 3) Per-student variant assignment (column L) -- only a subset of rows has
#    an assigned variant, the rest stay blank.
# ---------------------------------------------------------------------------
$variants = @{
  5  = 1
  9  = 1
  12 = 2
  16 = 3
  18 = 3
  23 = 1
  24 = 4
  28 = 3
  29 = 3
}

for ($r = 5; $r -le 33; $r++) {
    # Column K used to hold "=SUM(C{r}:I{r})" -- drop the formula/value, the
    # column becomes a leftover empty (still bordered/shaded) cell.
    $ws.Range("K$r").ClearContents()

    if ($variants.ContainsKey($r)) {
        $ws.Range("L$r").Value = $variants[$r]
    }

    # Column M..P now flag which variant (1..4) this row's student got.
    $ws.Range("M$r").Formula = "=IF(`$L$r=M`$3,1,0)"
    $ws.Range("N$r").Formula = "=IF(`$L$r=N`$3,1,0)"
    $ws.Range("O$r").Formula = "=IF(`$L$r=O`$3,1,0)"
    $ws.Range("P$r").Formula = "=IF(`$L$r=P`$3,1,0)"
}

# N4:P4 and M5:P33 take on the plain/no-border style used throughout the new
# helper columns.
$ws.Range("N4:P4").Font.Bold = $false
$ws.Range("M5:P33").Font.Bold = $false

# ---------------------------------------------------------------------------
# 4) Remove the old conditional formatting over K5:K33 and M5:M33 (it no
#    longer applies now that K is blank and M holds 0/1 flags).
# ---------------------------------------------------------------------------
$ws.Range("K5:K33").FormatConditions.Delete()
$ws.Range("M5:M33").FormatConditions.Delete()

# ---------------------------------------------------------------------------
# 5) View state: frozen pane moved down a few rows, and selection moved.
# ---------------------------------------------------------------------------
$ws.Range("J17").Select()
$excel.ActiveWindow.ScrollRow = 11

Write-Output "done"
